$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5
$ws1.Range("C2").Value = 0.5
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.6666666666666666
$ws1.Range("F2").Value = 0.8333333333333334
$ws1.Range("G2").Value = 0.9629629629629629
$ws1.Range("H2").Value = 0.7800905469286986
$ws1.Range("I2").Value = 534
$ws1.Range("J2").Value = 534
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0

$ws2.Range("B3").Value = 0.5
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.6666666666666666

$ws2.Range("B4").Value = 0.5
$ws2.Range("C4").Value = 0.5
$ws2.Range("D4").Value = 0.5
$ws2.Range("E4").Value = 0.5

$ws2.Range("B5").Value = 0.25
$ws2.Range("C5").Value = 0.5
$ws2.Range("D5").Value = 0.3333333333333333

$ws2.Range("B6").Value = 0.25
$ws2.Range("C6").Value = 0.5
$ws2.Range("D6").Value = 0.3333333333333333

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 534

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 534
